$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting existing rows 134:213 down to 135:214
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new weekly data point
$ws.Cells.Item(134,1).Value = 8
$ws.Cells.Item(134,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(134,3).Value = "Coquimbo"
$ws.Cells.Item(134,4).Value = 45001
$ws.Cells.Item(134,5).Value = 4
$ws.Cells.Item(134,6).Value = 100112040
$ws.Cells.Item(134,7).Value = "Cilantro"
$ws.Cells.Item(134,8).Value = "Sin especificar"
$ws.Cells.Item(134,9).Value = "Primera"
$ws.Cells.Item(134,10).Value = 2400
$ws.Cells.Item(134,11).Value = 2000
$ws.Cells.Item(134,12).Value = 2500
$ws.Cells.Item(134,13).Value = 2250
$ws.Cells.Item(134,14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(134,15).Value = "Provincia del Elquí"
$ws.Cells.Item(134,16).Value = 1500
$ws.Cells.Item(134,17).Value = 1.5
$ws.Cells.Item(134,18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Cells.Item(134,4).NumberFormat = $ws.Cells.Item(135,4).NumberFormat
